$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: updated to new publish date/time
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: now set to true (force text, not boolean)
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "true"
